# Sync automatico del tracker (cada 3h)
# Appends the latest batch of tracked picks to the bottom of Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=43; EventId="14580793"; Fecha="2025-09-01"; JugadorA="August Holmgren";    JugadorB="Thiago Agustin Tirante"; Pronostico="Gana August Holmgren";         Cuota=3.4 },
    @{ Row=44; EventId="14580346"; Fecha="2025-09-01"; JugadorA="Clement Chidekh";    JugadorB="Valentin Vacherot";       Pronostico="Gana Clement Chidekh";         Cuota=2.5 },
    @{ Row=45; EventId="14579620"; Fecha="2025-09-01"; JugadorA="Lukas Klein";        JugadorB="Saba Purtseladze";        Pronostico="Gana Saba Purtseladze";        Cuota=4 },
    @{ Row=46; EventId="14579621"; Fecha="2025-09-01"; JugadorA="Kaya Arinc";         JugadorB="Mert Naci Turker";        Pronostico="Gana Kaya Arinc";              Cuota=3.5 },
    @{ Row=47; EventId="14581354"; Fecha="2025-08-31"; JugadorA="Miguel Damas";       JugadorB="Szymon Kielan";           Pronostico="Gana Szymon Kielan";           Cuota=5 },
    @{ Row=48; EventId="14581358"; Fecha="2025-08-31"; JugadorA="Raul Brancaccio";    JugadorB="Ivan Marrero Curbelo";    Pronostico="Gana Ivan Marrero Curbelo";    Cuota=5.5 },
    @{ Row=49; EventId="14583788"; Fecha="2025-09-01"; JugadorA="Hiroki Moriya";      JugadorB="Yaojie Zeng";             Pronostico="Gana Yaojie Zeng";             Cuota=3.5 },
    @{ Row=50; EventId="14583789"; Fecha="2025-09-01"; JugadorA="Joshua Charlton";    JugadorB="Renta Tokuda";            Pronostico="Gana Joshua Charlton";         Cuota=3.4 },
    @{ Row=51; EventId="14583779"; Fecha="2025-09-01"; JugadorA="Tsung-Hao Huang";    JugadorB="Sergey Fomin";            Pronostico="Gana Tsung-Hao Huang";         Cuota=2.2 }
)

foreach ($r in $rows) {
    $row = $r.Row

    # event_id arrives from the API as a string id, keep it as text
    $cellA = $ws.Cells.Item($row, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $r.EventId

    # fecha is stored as plain text (yyyy-mm-dd), not a date serial
    $cellB = $ws.Cells.Item($row, 2)
    $cellB.NumberFormat = "@"
    $cellB.Value = $r.Fecha

    $ws.Cells.Item($row, 3).Value = $r.JugadorA
    $ws.Cells.Item($row, 4).Value = $r.JugadorB
    $ws.Cells.Item($row, 5).Value = $r.Pronostico
    $ws.Cells.Item($row, 6).Value = $r.Cuota

    # resultado / profit are pending until the match is settled
    $ws.Cells.Item($row, 7).Value = ""
    $ws.Cells.Item($row, 8).Value = ""
}
